$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column M ("national_id") -------------------------------------
# Copy column L's formatting (header + all 10 data rows, incl. borders)
# into column M so the new column matches the existing "status" column look.
$ws.Range("L1:L10").Copy()
$ws.Range("M1").PasteSpecial(-4122)  # xlPasteFormats

# Header
$ws.Range("M1").Value = "national_id"

# Row 2 gets a sample national id value, row 3 is intentionally left blank
# (matches the invalid-row test fixture semantics).
$ws.Range("M2").Value = "123"

# --- Column widths ----------------------------------------------------------
# D:E narrow slightly (23.4531 -> 23.5 char units)
$ws.Columns.Item(4).ColumnWidth = 22.63
$ws.Columns.Item(5).ColumnWidth = 22.63

# New column M should look like the K:L "wide" columns (~36.6719)
$ws.Columns.Item(13).ColumnWidth = 35.8

Write-Output "done"
